$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the three tables (slides 14, 15, 16) to the built-in table
#    style {BD2F32B0-2BE6-4A50-8853-602F58A1E120} (was
#    {CCB64940-EC2D-4855-AB3C-D287CA2A6AE7}).
# ---------------------------------------------------------------------------
$newTableStyle = "{BD2F32B0-2BE6-4A50-8853-602F58A1E120}"
foreach ($slideIdx in 14,15,16) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newTableStyle)
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the two themes used by the deck: the slide master (theme1.xml)
#    was the "Integral / Red Violet" palette and becomes the default
#    "Office" palette that used to live only on the notes master
#    (theme2.xml).
# ---------------------------------------------------------------------------
function Set-ThemeColor {
    param($scheme, [int]$index, [string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$slideForTheme = $p.Slides.Item(1)
$themeColors = $slideForTheme.ThemeColorScheme

Set-ThemeColor $themeColors 1  "000000"
Set-ThemeColor $themeColors 2  "FFFFFF"
Set-ThemeColor $themeColors 3  "44546A"
Set-ThemeColor $themeColors 4  "E7E6E6"
Set-ThemeColor $themeColors 5  "5B9BD5"
Set-ThemeColor $themeColors 6  "ED7D31"
Set-ThemeColor $themeColors 7  "A5A5A5"
Set-ThemeColor $themeColors 8  "FFC000"
Set-ThemeColor $themeColors 9  "4472C4"
Set-ThemeColor $themeColors 10 "70AD47"
Set-ThemeColor $themeColors 11 "0563C1"
Set-ThemeColor $themeColors 12 "954F72"
